# [Kadastro App] Yeni kayit eklendi: 2958
$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")
$newRowValues = @("2958", "2025-09-09", "Erdemli", "1", "ÇAP", "CEMAL TİMUROĞLU (K.Teknisyeni)")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 49
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($newRow, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $newRowValues[$col - 1]
        $cell.Style = "Normal"
    }
}
